$p = $ppt.ActivePresentation

# Slides 4 and 5 both contain a "Content Placeholder" text box (shape index 2)
# whose text repeats a "<side>ty texty izquierdi texti" pattern. Every
# "izquierdi" run is a typo that should read "derechi".
foreach ($slideIdx in 4, 5) {
    $s = $p.Slides.Item($slideIdx)
    $sh = $s.Shapes.Item(2)
    $tr = $sh.TextFrame.TextRange
    $tr.Replace("izquierdi", "derechi")
}

# Slide 6: remove the empty, unused "Content Placeholder 5" textbox (the one
# with the creationId B9EDF57C-D77E-5F46-B22C-762045F7D38F) that sits between
# the "Content Placeholder 7" box and the final "Content Placeholder 5" box.
$s6 = $p.Slides.Item(6)
for ($i = $s6.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s6.Shapes.Item($i)
    if ($sh.Name -eq "Content Placeholder 5" -and $sh.Id -eq 10) {
        $sh.Delete()
    }
}
